$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E, F, L, M, O across rows 2-25
# as per diff (Case_5_246, "case with 380 kV done")
$cols = @("B","C","D","E","F","L","M","O")

$data = @(
    @(15.8291344720034, 11.81307979361698, 5.239497140747618, 12.73574010039477, 25.99059951873788, 9.876869094851983, 14.88139831220169, 23.25365463924795),
    @(15.28223480662598, 11.61793548042829, 5.208347467135989, 12.78793604549683, 26.04128224526769, 9.884781208495161, 14.76221821450033, 23.34896619047704),
    @(14.93775673891448, 11.49608816262041, 5.189065526473325, 12.82165585133474, 26.08123352351211, 9.890994933423983, 14.69048682200406, 23.41422011713684),
    @(14.79540284180102, 11.4459684698061, 5.181172078059258, 12.83581845923576, 26.09972578195558, 9.89386838780117, 14.66164327331016, 23.44249747632663),
    @(14.77165212672973, 11.43761930416882, 5.179859304918836, 12.83819564477836, 26.10292970653478, 9.894366148080307, 14.65687792274414, 23.4472945224988),
    @(14.93584461976241, 11.49541405830676, 5.188959213685112, 12.82184514486052, 26.08147397399739, 9.891032303371263, 14.69009622725584, 23.4145946589983),
    @(15.64248268655192, 11.74624091424196, 5.228791076715746, 12.75339114463158, 26.00623676707072, 9.879316095016046, 14.84001958004105, 23.28511662115811),
    @(16.9512377887164, 12.22016117634855, 5.305519951005589, 12.6323564086523, 25.92912536055191, 9.86707704769532, 15.14440715931313, 23.08493323893319),
    @(17.85633472017371, 12.55503191110265, 5.360859205103844, 12.55140038061511, 25.91580649815194, 9.864600721555158, 15.37289493485608, 22.97101314357257),
    @(18.25416794598627, 12.70401346838582, 5.385769527594006, 12.51628437560353, 25.91921473818765, 9.86488134335559, 15.47757343641451, 22.92646853680101),
    @(18.40270121067948, 12.75991017816015, 5.395161043943726, 12.50323164513862, 25.92186908058317, 9.865189154930462, 15.51729146619905, 22.91065351360698),
    @(18.37080782884123, 12.74789549217636, 5.393140314043099, 12.50603190859919, 25.92123674699777, 9.86511391186591, 15.50873442301023, 22.91401262766512),
    @(18.26643093408272, 12.70862273943036, 5.386543024905286, 12.51520561719948, 25.91940577337897, 9.864902632743055, 15.48083964620831, 22.92514628475459),
    @(18.20221803122676, 12.68449833613771, 5.382496487677887, 12.52085664303384, 25.91846189004604, 9.864799440521773, 15.46376270157293, 22.93210329582382),
    @(17.8300451920178, 12.54522491896428, 5.359225644894199, 12.55372963180693, 25.91577449901384, 9.864610630798035, 15.36606644330323, 22.97407125364099),
    @(17.59807987380638, 12.45889931323018, 5.344879641864949, 12.57433363314822, 25.91655276404789, 9.86485465676596, 15.30630238751051, 23.00168605410821),
    @(17.46335734665409, 12.40893362511766, 5.336603582992225, 12.58634565303552, 25.91789143738125, 9.865127510826955, 15.2719988254647, 23.01825405038189),
    @(17.41752281989676, 12.39196342849759, 5.333797332199671, 12.59044042810787, 25.91849762919358, 9.865242673939063, 15.26039728935978, 23.02398108422641),
    @(17.62290865891207, 12.4681215418308, 5.34640937197404, 12.57212363262107, 25.91637767566769, 9.864814972381138, 15.31265721047108, 22.99867550046244),
    @(18.29714730171769, 12.72017248123893, 5.388481962625083, 12.5125044381196, 25.91990655238268, 9.864959227358318, 15.489031106119, 22.92184743039018),
    @(18.72541353469959, 12.8818592054338, 5.415735415948861, 12.47496701296436, 25.93016169913691, 9.866227848147266, 15.60474677543784, 22.87777660556677),
    @(18.49800877429047, 12.79585418151813, 5.401213173570462, 12.49487123309995, 25.92396062076321, 9.865443591728202, 15.54295551341073, 22.90073410582023),
    @(17.61168780012356, 12.46395321586088, 5.345717869014935, 12.57312225511589, 25.91645405713726, 9.864832500720432, 15.30978401974065, 23.00003441705061),
    @(16.60647758000777, 12.09412687057498, 5.284930391728268, 12.66369440177342, 25.94239966343726, 9.86924138674353, 15.06110139386856, 23.1332940895555)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $cellAddr = "$($cols[$j])$row"
        $ws.Range($cellAddr).Value = $rowValues[$j]
    }
}
